$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert three new rows at 12-14 (pushes the "Programa resumido:" block and
# everything below it down by three rows), to make room for a new
# "Docentes responsáveis:" entry with two names.
$ws.Rows("12:14").Insert()

# Row 12: label only (mirrors the "Avaliação:"/"Requisitos:" label-only rows)
$ws.Range("A12").Value = "Docentes responsáveis:"

# Row 13: first professor, duplicated into B and C like the rest of the sheet
$ws.Range("B13").Value = "3480026 - João Paulo Pascon"
$ws.Range("C13").Value = "3480026 - João Paulo Pascon"

# Row 14: second professor, duplicated into B and C
$ws.Range("B14").Value = "1176388 - Luiz Tadeu Fernandes Eleno"
$ws.Range("C14").Value = "1176388 - Luiz Tadeu Fernandes Eleno"

# The row insert stamps every styled column (A/B/C) with a blank, formatted
# cell. Clear the cells that should stay completely empty so the row only
# contains the cells actually used (matching the rest of the sheet's
# label-only / data-only row shapes).
$ws.Range("B12").Clear()
$ws.Range("C12").Clear()
$ws.Range("A13").Clear()
$ws.Range("A14").Clear()
